$d = $word.ActiveDocument

# 1. Delete the empty centered paragraph (Times New Roman, size 36) that
#    sits just before the "Logistic Regression" Heading1 paragraph near the
#    top of the document.
$p = $d.Paragraphs(24)
$p.Range.Delete()

# 2. Merge "depended variable using a given set of " + "independent" +
#    " variable." into a single run "depended variable using a given set of
#    independent variable." (text is unchanged, only run/proofErr structure
#    is simplified) without touching the preceding run.
$rng = $d.Content
[void]$rng.Find.Execute("depended variable using a given set of independent variable.")
$s = $rng.Start
$e = $rng.End
$target = $d.Range($s, $e)
$target.Delete()
$insertPoint = $d.Range($s, $s)
$insertPoint.InsertAfter("depended variable using a given set of independent variable.")

# 3. Find "The main thing is..." paragraph (about raining probability) and
#    insert a new blank paragraph right after it.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    if ($pp.Range.Text -match "raining today") {
        $idx = $i
        break
    }
}
$d.Paragraphs($idx).Range.InsertParagraphAfter()

# 4. The next two paragraphs (previously plain empty paragraphs) get a
#    yellow highlight on their paragraph mark, matching the highlighted
#    paragraphs around them.
$d.Paragraphs($idx + 2).Range.HighlightColorIndex = 7
$d.Paragraphs($idx + 3).Range.HighlightColorIndex = 7
